$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# --- Main document body (Simplified -> Traditional Chinese) ---
Replace-Text "向 ROW 客户发送提醒电子邮件" "向 ROW 客戶傳送提醒電子郵件"
Replace-Text "将于 9 月 29 日移除 Tether Omni (USDT)" "將於 9 月 29 日移除 Tether Omni (USDT)"
Replace-Text "向 Tether Omni 说再见" "向 Tether Omni 道別"
Replace-Text "自 2023 年 9 月 29 日格林威治标准时间 00:00 起，Deriv 将停止提供Tether Omni (USDT) 作为账户货币。 这是因为 Tether 已停止支持 USDT 的 Omni 转账。" "自 2023 年 9 月 29 日格林威治標準時間 00:00 起，Deriv 將停止提供Tether Omni (USDT) 作為帳戶貨幣。 這是因為 Tether 已停止支援 USDT 的 Omni 轉帳。"
Replace-Text "需要做什么？" "需要做什麼？"
Replace-Text "如果 USDT 账户 " "如果 USDT 帳戶 "
Replace-Text " 中有余额，请在上述日期之前提取余额。 如果有持仓头寸，提取余额之前请先平仓。" " 中有餘額，請在上述日期之前提取餘額。 若有持倉頭寸，提取餘額前請先平倉。"
Replace-Text "查看账户" "檢查帳戶"
Replace-Text "USDT 账户将于 2023 年 9 月 29 日格林尼治标准时间 00:00 关闭。 任何持仓头寸将在上述日期后自动平仓，账户余额将转移到最后活跃的账户" "USDT 帳戶將於 2023 年 9 月 29 日格林尼治標準時間 00:00 關閉。 任何持倉頭寸將在上述日期後自動平倉，帳戶餘額將轉移到最後活躍的帳戶"
Replace-Text "在此过程中将采用标准汇率和费用。" "在此過程中將採用標準匯率和費用。"
Replace-Text "如有任何疑问，请通过以下方式联系我们：" "如有任何疑問，請透過以下方式聯繫我們："
Replace-Text "实时聊天" "即時聊天"

# --- Comments (Simplified -> Traditional Chinese) ---
# Comment bodies live in a separate story that this host's Find/Range
# machinery does not expose through $d.Content, so walk the Comments
# collection directly and update each comment's Range text in place.
$commentEdits = @(
    @{ Old = "@azita@regentmarkets.com，BE 不能保证届时可以完成脚本，"; New = "@azita@regentmarkets.com，BE 無法保證屆時可以完成指令，" },
    @{ Old = "可以将其称为“所述日期之后”吗？"; New = "可以將其稱為“所述日期之後”嗎？" },
    @{ Old = "意思是在所述日期之后完成转账吗？"; New = "意思是在所述日期之後完成轉帳嗎？" },
    @{ Old = "是的... 看来我们无法确认日期"; New = "是的... 看來我們無法確認日期" }
)

for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $c = $d.Comments.Item($i)
    foreach ($edit in $commentEdits) {
        try {
            if ($c.Range.Text -eq $edit.Old) {
                $c.Range.Text = $edit.New
            }
        } catch {
        }
        try {
            $null = $c.Range.Find.Execute($edit.Old, $true, $false, $false, $false, $false, $true, 1, $false, $edit.New, 2)
        } catch {
        }
    }
}
